# Add a new "VWC_perc" column (C) to each species sheet: header in C1,
# and C2:C61 = corresponding B value / 10 (percent volumetric water content).
# Also drop the two duplicate "_xlchart.v1.8" / "_xlchart.v1.9" defined
# names (exact duplicates of .6 / .7), and leave the selection/active-sheet
# state the way it was left after making these edits (PIPO tab active with
# the new column selected; other sheets' selection moved to the new column).

$wb = $excel.ActiveWorkbook

$sheetNames = @("PIPO", "PSME", "PIFL", "PIEN")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value = "VWC_perc"
    $ws.Range("C2").Formula = "=B2/10"
    $ws.Range("C3:C61").Formula = "=B3/10"
}

# Remove the two stray duplicate defined names pointing at PIEN!$B$1 /
# PIEN!$B$2:$B$61 (exact dupes of _xlchart.v1.6 / _xlchart.v1.7).
$wb.Names.Item("_xlchart.v1.8").Delete()
$wb.Names.Item("_xlchart.v1.9").Delete()

$ws1 = $wb.Worksheets.Item("PIPO")
$ws2 = $wb.Worksheets.Item("PSME")
$ws3 = $wb.Worksheets.Item("PIFL")
$ws4 = $wb.Worksheets.Item("PIEN")

# Walk through the sheets leaving each one's selection on/near the new
# column, then land back on PIPO as the active tab.
$ws4.Activate()
$ws4.Range("D18").Select()
$win = $excel.ActiveWindow
$win.Zoom = 123

$ws3.Activate()
$ws3.Range("D8").Select()

$ws2.Activate()
$ws2.Range("D2").Select()

$ws1.Activate()
$ws1.Range("C2:C61").Select()
